# The deck ships with two theme parts:
#   ppt/theme/theme1.xml -> bound to the (only) slide master, currently the
#                            "Integral" / "Red Violet" look
#   ppt/theme/theme2.xml -> bound to the notes master, currently the plain
#                            default "Office Theme" look
#
# The target edit swaps the two palettes: the slide master should switch to
# the default blue "Office" color scheme, while the notes master would take
# on the old "Red Violet" colors. PowerPoint's automation model doesn't give
# us a raw-XML/theme-import hook, so we reproduce the swap using the real
# ThemeColorScheme object (the 12-slot modern DrawingML color scheme -
# dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink, in that order) exposed per slide.

function ConvertTo-BgrColor([string]$hex) {
    # PowerPoint's RGB color properties take a BGR-packed long (standard
    # VBA RGB() encoding), so pack the hex RRGGBB triplet accordingly.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target palette: the stock Office theme colors, in clrScheme order.
$officePalette = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation

# All slides share the single slide master/theme, so editing the scheme via
# any slide updates the master's underlying theme part (theme1.xml).
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officePalette.Count; $i++) {
    $themeColor = $themeColors.Colors($i)
    $themeColor.RGB = ConvertTo-BgrColor $officePalette[$i - 1]
}
